# "make the battle scene can enter"
# The MainIcon sheet's Icon/path column (N) used PascalCase resource names
# (MainIconX / SideButton4). The game expects lower-case path strings, so
# normalize them to lower-case (mainiconX / sidebutton4) so the referenced
# scene assets resolve and the battle scene can be entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N4").Value  = "mainicon1"
$ws.Range("N5").Value  = "mainicon8"
$ws.Range("N6").Value  = "mainicon3"
$ws.Range("N7").Value  = "mainicon5"
$ws.Range("N8").Value  = "mainicon2"
$ws.Range("N9").Value  = "mainicon10"
$ws.Range("N10").Value = "sidebutton4"

# Move the active selection the way the authored workbook ended up (N11).
$ws.Range("N11").Select()
